$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.030943751335144
$ws.Range("B1").Value = 2.293096542358398
$ws.Range("C1").Value = 9.679383277893066
$ws.Range("D1").Value = 2.46974778175354
$ws.Range("E1").Value = 1.357510328292847
